# Apply the commit "Creacion inicial de los modelos del reto mas el modelo Empresas Cliente"
# to Organizacion.xlsx (sheet "Hoja1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: widen it to fit the (now longer) task descriptions ------------------
$ws.Columns.Item(2).ColumnWidth = 44.7109375

# --- Row 5 (task #4): replace the old "Disenar modelos" task with the new one ------
$ws.Range("B5").Value = "Creacion inicial de los modelos del reto con la adicion de Clientes"
$ws.Range("C5").Value = "Jon"
$ws.Range("D5").Value = 45778
$ws.Range("E5").Value = 45779
$ws.Range("F5").Value = "✅ Hecho"
$ws.Range("G5").Value = ""
$ws.Rows.Item(5).RowHeight = 45

# --- Row 6 (task #5): comment cell style simplified, row height back to default ----
# (copy the formatting that G7 below will need from an existing "style 4" cell so we
# do not introduce brand-new, unused cell formats into the workbook)
$ws.Range("F5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Rows.Item(6).AutoFit()

# --- Row 7 (task #6): fill in the new task + comment -------------------------------
$ws.Range("B7").Value = "Crear vistas del detalle de cada modelo y cambiar el link en el headder de la plantilla base para navegar entre ellos "
$ws.Range("F5").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = "Probablemente hacer cada uno uno"
$ws.Rows.Item(7).RowHeight = 45

$excel.CutCopyMode = 0

# --- Selection cosmetics -------------------------------------------------------------
$null = $ws.Range("G24").Select()
